$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style from an existing header cell (AB1) onto the new
# header cells so they pick up the exact same cellXf (bold, bordered,
# centered) instead of Excel synthesizing a near-duplicate style.
$ws.Range("AB1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Season record (Wins/Losses/Ties) is the same for every player row on
# this sheet (it's the team's overall record for the season), so fill
# AC2:AE37 row by row.
for ($r = 2; $r -le 37; $r++) {
    $ws.Cells.Item($r, 29).Value = 49
    $ws.Cells.Item($r, 30).Value = 64
    $ws.Cells.Item($r, 31).Value = 0
}
